# Updated cryptos list — refresh Price (D) and Volume(1h) (E) columns.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $text) {
    # Force the cell to stay a text cell (many of these "prices" are
    # dotted strings like 41.910.60 that must not be parsed as numbers).
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.ClearFormats()
}

Set-TextValue $ws.Range("D2") "41.910.60"
Set-TextValue $ws.Range("D3") "2.285.21"
$ws.Range("E3").Value = "  -3.08%  "
$ws.Range("E4").Value = "  -0.07%  "
Set-TextValue $ws.Range("D5") "316.39"
$ws.Range("E5").Value = "  -0.73%  "
Set-TextValue $ws.Range("D6") "103.16"
$ws.Range("E6").Value = "  -3.44%  "
$ws.Range("E7").Value = "  -0.77%  "
$ws.Range("E8").Value = "  +0.05%  "
Set-TextValue $ws.Range("D9") "0.602"
$ws.Range("E9").Value = "  -2.74%  "
Set-TextValue $ws.Range("D10") "39.22"
$ws.Range("E10").Value = "  -4.90%  "
Set-TextValue $ws.Range("D11") "0.0897"
$ws.Range("E11").Value = "  -3.00%  "
Set-TextValue $ws.Range("D12") "8.22"
$ws.Range("E12").Value = "  -3.59%  "
$ws.Range("E13").Value = "  -0.37%  "
Set-TextValue $ws.Range("D14") "0.956"
$ws.Range("E14").Value = "  -4.32%  "
Set-TextValue $ws.Range("D15") "15.13"
$ws.Range("E15").Value = "  -5.46%  "
Set-TextValue $ws.Range("D16") "2.632.49"
$ws.Range("E16").Value = "  -3.17%  "
Set-TextValue $ws.Range("D17") "2.293.65"
$ws.Range("E17").Value = "  -3.23%  "
Set-TextValue $ws.Range("D18") "42.007.00"
Set-TextValue $ws.Range("D19") "7.38"
$ws.Range("E19").Value = "  -2.49%  "
$ws.Range("E20").Value = "  -0.98%  "
Set-TextValue $ws.Range("D21") "3.62"
$ws.Range("E21").Value = "  -0.98%  "
Set-TextValue $ws.Range("D22") "73.15"
$ws.Range("E22").Value = "  -3.80%  "
Set-TextValue $ws.Range("D23") "278.46"
$ws.Range("E23").Value = "  +3.63%  "
Set-TextValue $ws.Range("D24") "10.06"
$ws.Range("E24").Value = "  +6.73%  "
Set-TextValue $ws.Range("D25") "2.25"
$ws.Range("E25").Value = "  -2.68%  "
$ws.Range("E26").Value = "  +0.83%  "
Set-TextValue $ws.Range("D27") "10.75"
$ws.Range("E27").Value = "  -5.74%  "
$ws.Range("E28").Value = "  +4.33%  "
Set-TextValue $ws.Range("D29") "22.76"
$ws.Range("E29").Value = "  -2.98%  "
Set-TextValue $ws.Range("D30") "35.71"
$ws.Range("E30").Value = "  -3.24%  "
Set-TextValue $ws.Range("D31") "162.63"
$ws.Range("E31").Value = "  -3.04%  "
$ws.Range("E32").Value = "  -4.10%  "
$ws.Range("E33").Value = "  -1.91%  "
Set-TextValue $ws.Range("D34") "5.81"
$ws.Range("E34").Value = "  -2.75%  "
$ws.Range("E35").Value = "  +3.75%  "
$ws.Range("E36").Value = "  -6.35%  "
Set-TextValue $ws.Range("D37") "4.49"
$ws.Range("E37").Value = "  -5.15%  "
Set-TextValue $ws.Range("D38") "0.0346"
$ws.Range("E38").Value = "  -4.51%  "
$ws.Range("E39").Value = "  +3.85%  "
Set-TextValue $ws.Range("D40") "3.73"
$ws.Range("E40").Value = "  -2.53%  "
Set-TextValue $ws.Range("D41") "98.87"
$ws.Range("E41").Value = "  -7.61%  "
$ws.Range("E42").Value = "  -4.64%  "
Set-TextValue $ws.Range("D43") "68.99"
$ws.Range("E43").Value = "  -2.97%  "
$ws.Range("E44").Value = "  +0.10%  "
$ws.Range("E45").Value = "  -6.19%  "
Set-TextValue $ws.Range("D46") "11.86"
$ws.Range("E46").Value = "  -3.71%  "
Set-TextValue $ws.Range("D47") "112.08"
$ws.Range("E47").Value = "  -1.01%  "
Set-TextValue $ws.Range("D48") "76.78"
$ws.Range("E48").Value = "  +1.52%  "
Set-TextValue $ws.Range("D49") "8.88"
$ws.Range("E49").Value = "  -2.87%  "
Set-TextValue $ws.Range("D50") "5.26"
$ws.Range("E50").Value = "  -4.81%  "
Set-TextValue $ws.Range("D51") "1.575.07"
$ws.Range("E51").Value = "  -0.02%  "
